$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 4 new daily rows (18-21 Apr 2021 / serials 44304-44307) below the
# existing data, mirroring the layout/style of the last existing row (229).
$dates = @(44304, 44305, 44306, 44307)
$bvals = @(1, 0, 0, 0)
$cvals = @(9, 9, 9, 9)
$dvals = @(274.8091603053435, 274.8091603053435, 274.8091603053435, 274.8091603053435)

$lastRow = 229

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $lastRow + 1 + $i

    # Reuse the date cell's style (border/alignment/bold + date number format)
    # from the row above, so the new date cells keep style index 2, just like
    # every other date cell in column A.
    $ws.Cells.Item($lastRow, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $bvals[$i]
    $ws.Cells.Item($r, 3).Value = $cvals[$i]
    $ws.Cells.Item($r, 4).Value = $dvals[$i]
}
